# Swap the order of the last two slides in the deck:
#   - the "Generalized Polynomial Chaos" slide (currently slide 29)
#   - the "raven.gif" picture slide (currently slide 30)
# After this edit, the raven.gif picture slide should be slide 29 and the
# Generalized Polynomial Chaos slide should be slide 30.

$p = $ppt.ActivePresentation

# Move what is currently the last slide (30, raven.gif) to be slide 29,
# pushing the current slide 29 (Generalized Polynomial Chaos) down to 30.
$s = $p.Slides.Item(30)
$s.MoveTo(29)
